$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 4 new blank columns so the sheet grows from A:L to A:P.
# New columns must land (in the FINAL layout) at F, J, M and P.
# Inserting left-to-right using the FINAL column letters works because each
# insert only shifts cells to the right of its own insertion point, and we
# have not yet touched anything at/after the next (further-right) target.
# ---------------------------------------------------------------------------
$ws.Columns("F").Insert() | Out-Null
$ws.Columns("J").Insert() | Out-Null
$ws.Columns("M").Insert() | Out-Null
$ws.Columns("P").Insert() | Out-Null

# ---------------------------------------------------------------------------
# Give the 4 new header cells the same look (bold / border / centered) as the
# rest of row 1 by copying an existing header cell's formatting over them.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Header row text (new columns only -- the others were already shifted into
# place correctly by the column inserts above).
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "rmse_vcell_hybrid_without_delta [mV]"
$ws.Range("J1").Value = "rmse_thetass2_hybrid_without_delta [milli]"
$ws.Range("M1").Value = "rmse_phie2_hybrid_without_delta [mV]"
$ws.Range("P1").Value = "rmse_if2_hybrid_without_delta [mA]"

# ---------------------------------------------------------------------------
# Data rows 2..8.
#   E = rmse_vcell_hybrid        (value changes)
#   F = rmse_vcell_hybrid_without_delta   (new column)
#   I = rmse_thetass2_hybrid     (value changes)
#   J = rmse_thetass2_hybrid_without_delta (new column)
#   L = rmse_phie2_hybrid        (value changes)
#   M = rmse_phie2_hybrid_without_delta   (new column)
#   O = rmse_if2_hybrid          (value changes)
#   P = rmse_if2_hybrid_without_delta     (new column)
# ---------------------------------------------------------------------------
$data = @{
    2 = @{ E=0.2484402900329776; F=1.41983471626149;  I=0.267014520436426;  J=1.996485785104415; L=0.1915680168484054; M=0.9584283649439964; O=1.986758250353938;  P=4.853599605237231 }
    3 = @{ E=0.1992275861903888; F=1.683943637335505;  I=0.2291191651059691; J=1.50508416648509;  L=0.1094559105373297; M=0.2673031981268261; O=0.5981098243075972; P=12.18018861377182 }
    4 = @{ E=0.2957377338123012; F=2.272985895707731;  I=0.3188370031423483; J=2.969493742516178; L=0.1195641350490709; M=1.362948717325113;  O=1.005707166316514;  P=21.53801454238877 }
    5 = @{ E=4.062989119600359;  F=7.638299160277519;  I=3.194834529683527;  J=5.731017690657275; L=0.2663234867236759; M=5.830591755813333;  O=3.037047053381891;  P=11.25843421902136 }
    6 = @{ E=0.3083678849325592; F=1.804530396135378;  I=0.652767866820271;  J=13.92147349229983; L=0.2568417415037383; M=2.235569298748505;  O=1.37121513926494;   P=47.73187631198127 }
    7 = @{ E=0.5170291484078821; F=3.376109704298808;  I=0.6092160685244821; J=12.07311562416297; L=0.2042234060781455; M=2.529502188984647;  O=3.334700528819338;  P=47.47747159958216 }
    8 = @{ E=1.694119321374044;  F=6.217035878164081;  I=1.478611419389102;  J=14.15199215348284; L=0.7096032535990036; M=2.406538298692395;  O=10.79875824800481;  P=38.58328398286862 }
}

foreach ($r in 2..8) {
    $row = $data[$r]
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
}
